$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.2484426666666667
$ws.Range("H2").Value = 0.745328
$ws.Range("I2").Value = 0.1396403772415532
$ws.Range("J2").Value = 0.1396403772415532
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.498163333333333
$ws.Range("N2").Value = 4.49449
$ws.Range("O2").Value = 0.02101839619520399
$ws.Range("P2").Value = 0.021018396195204
$ws.Range("Q2").Value = 0.3722076936355556
$ws.Range("R2").Value = 3.34986924272
$ws.Range("S2").Value = 0.002935016773710713
$ws.Range("T2").Value = 0.002935016773710713

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.2484426666666667
$ws.Range("H3").Value = 0.745328
$ws.Range("I3").Value = 0.1396403772415532
$ws.Range("J3").Value = 0.1396403772415532
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.96588933333333
$ws.Range("N3").Value = 83.897668
$ws.Range("O3").Value = 0.3923458336491322
$ws.Range("P3").Value = 0.3923458336491322
$ws.Range("Q3").Value = 6.947920121678222
$ws.Range("R3").Value = 62.531281095104
$ws.Range("S3").Value = 0.05478732021991651
$ws.Range("T3").Value = 0.05478732021991651

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.2484426666666667
$ws.Range("H4").Value = 0.745328
$ws.Range("I4").Value = 0.1396403772415532
$ws.Range("J4").Value = 0.1396403772415532
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 41.81461766666666
$ws.Range("N4").Value = 125.443853
$ws.Range("O4").Value = 0.5866357701556637
$ws.Range("P4").Value = 0.5866357701556638
$ws.Range("Q4").Value = 10.38853511875378
$ws.Range("R4").Value = 93.49681606878399
$ws.Range("S4").Value = 0.081918040247926
$ws.Range("T4").Value = 0.08191804024792601

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.055305333333333
$ws.Range("H5").Value = 3.165916
$ws.Range("I5").Value = 0.5931478551122046
$ws.Range("J5").Value = 0.5931478551122047
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.498163333333333
$ws.Range("N5").Value = 4.49449
$ws.Range("O5").Value = 0.02101839619520399
$ws.Range("P5").Value = 0.021018396195204
$ws.Range("Q5").Value = 1.581019755871111
$ws.Range("R5").Value = 14.22917780284
$ws.Range("S5").Value = 0.01246701662108377
$ws.Range("T5").Value = 0.01246701662108378

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.055305333333333
$ws.Range("H6").Value = 3.165916
$ws.Range("I6").Value = 0.5931478551122046
$ws.Range("J6").Value = 0.5931478551122047
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.96588933333333
$ws.Range("N6").Value = 83.897668
$ws.Range("O6").Value = 0.3923458336491322
$ws.Range("P6").Value = 0.3923458336491322
$ws.Range("Q6").Value = 29.51255216487644
$ws.Range("R6").Value = 265.612969483888
$ws.Range("S6").Value = 0.2327190896911926
$ws.Range("T6").Value = 0.2327190896911926

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.055305333333333
$ws.Range("H7").Value = 3.165916
$ws.Range("I7").Value = 0.5931478551122046
$ws.Range("J7").Value = 0.5931478551122047
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 41.81461766666666
$ws.Range("N7").Value = 125.443853
$ws.Range("O7").Value = 0.5866357701556637
$ws.Range("P7").Value = 0.5866357701556638
$ws.Range("Q7").Value = 44.12718903492755
$ws.Range("R7").Value = 397.144701314348
$ws.Range("S7").Value = 0.3479617487999281
$ws.Range("T7").Value = 0.3479617487999283

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4754126666666667
$ws.Range("H8").Value = 1.426238
$ws.Range("I8").Value = 0.2672117676462422
$ws.Range("J8").Value = 0.2672117676462422
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.498163333333333
$ws.Range("N8").Value = 4.49449
$ws.Range("O8").Value = 0.02101839619520399
$ws.Range("P8").Value = 0.021018396195204
$ws.Range("Q8").Value = 0.7122458254022223
$ws.Range("R8").Value = 6.41021242862
$ws.Range("S8").Value = 0.005616362800409511
$ws.Range("T8").Value = 0.005616362800409512

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4754126666666667
$ws.Range("H9").Value = 1.426238
$ws.Range("I9").Value = 0.2672117676462422
$ws.Range("J9").Value = 0.2672117676462422
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.96588933333333
$ws.Range("N9").Value = 83.897668
$ws.Range("O9").Value = 0.3923458336491322
$ws.Range("P9").Value = 0.3923458336491322
$ws.Range("Q9").Value = 13.29533802366489
$ws.Range("R9").Value = 119.658042212984
$ws.Range("S9").Value = 0.1048394237380231
$ws.Range("T9").Value = 0.1048394237380231

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4754126666666667
$ws.Range("H10").Value = 1.426238
$ws.Range("I10").Value = 0.2672117676462422
$ws.Range("J10").Value = 0.2672117676462422
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 41.81461766666666
$ws.Range("N10").Value = 125.443853
$ws.Range("O10").Value = 0.5866357701556637
$ws.Range("P10").Value = 0.5866357701556638
$ws.Range("Q10").Value = 19.87919889055711
$ws.Range("R10").Value = 178.912790015014
$ws.Range("S10").Value = 0.1567559811078096
$ws.Range("T10").Value = 0.1567559811078096

